$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.706.30'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '2.557.26'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.62'
$ws.Range("E5").Value = '  +2.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.49'
$ws.Range("E6").Value = '  +3.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.573'
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.545'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.41'
$ws.Range("E10").Value = '  +2.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0805'
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.117'
$ws.Range("E12").Value = '  +9.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.51'
$ws.Range("E13").Value = '  -1.86%  '
$ws.Range("D14").Value = '2.572.25'
$ws.Range("E14").Value = '  +1.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.877'
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.32'
$ws.Range("E16").Value = '  +1.69%  '
$ws.Range("D17").Value = '42.763.57'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.18'
$ws.Range("E18").Value = '  +7.06%  '
$ws.Range("D19").Value = '0.0₃0983'
$ws.Range("E19").Value = '  +2.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.59'
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.41'
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '256.09'
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.94'
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.09'
$ws.Range("E24").Value = '  -1.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '27.96'
$ws.Range("E25").Value = '  -4.63%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.90'
$ws.Range("E27").Value = '  +7.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.97'
$ws.Range("E30").Value = '  +1.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.39'
$ws.Range("E31").Value = '  +4.29%  '
$ws.Range("B32").Value = 'WEMIXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.76'
$ws.Range("E32").Value = '  +2.04%  '
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.14'
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.31'
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0796'
$ws.Range("E35").Value = '  +0.73%  '
$ws.Range("B36").Value = 'EnergySwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '25.88'
$ws.Range("E36").Value = '  +7.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.114'
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.00'
$ws.Range("E38").Value = '  +13.02%  '
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.85'
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("B41").Value = 'ApeXProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.05'
$ws.Range("E41").Value = '  +29.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.36'
$ws.Range("E42").Value = '  -1.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0306'
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("D44").Value = '2.063.23'
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '87.68'
$ws.Range("E46").Value = '  +3.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.23'
$ws.Range("E47").Value = '  +5.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.63'
$ws.Range("E48").Value = '  +10.57%  '
$ws.Range("D49").Value = '2.812.30'
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.89'
$ws.Range("E50").Value = '  +1.18%  '
$ws.Range("E51").Value = '  +2.91%  '
